# "Adding papers back to 2011" -- restructure ConferenceStats.xlsx:
#   * Rename Sheet1 -> "By Topic"
#   * Add a totals row (row 13) to "By Topic"
#   * Add a brand new "By Conf and Year" sheet with a by-year/by-conference
#     breakdown (with its own Totals row/column), make it the active sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename the existing sheet and add the totals row to it.
# ---------------------------------------------------------------------
$byTopic = $wb.Worksheets.Item(1)
$byTopic.Name = "By Topic"

$byTopic.Range("A13").Value = "Totals"
$byTopic.Range("B13").Formula = "=SUM(B2:B12)"
$byTopic.Range("C13:L13").Formula = "=SUM(C2:C12)"

# ---------------------------------------------------------------------
# 2) Add the new sheet right after "By Topic" and fill it in.
# ---------------------------------------------------------------------
$byConfYear = $wb.Worksheets.Add($null, $byTopic)
$byConfYear.Name = "By Conf and Year"

# Header row: years across the top, "Totals" label in H1 (bold).
$byConfYear.Range("B1").Value = 2012
$byConfYear.Range("C1").Value = 2013
$byConfYear.Range("D1").Value = 2014
$byConfYear.Range("E1").Value = 2015
$byConfYear.Range("F1").Value = 2016
$byConfYear.Range("G1").Value = 2017
$byConfYear.Range("H1").Value = "Totals"
$byConfYear.Range("H1").Font.Bold = $true

# Per-conference rows of paper counts by year.
$confRows = @(
    @{ Row = 2;  Name = "CCS";       Vals = @(1, 0, 3, 8, 10, 9) },
    @{ Row = 3;  Name = "Oakland";   Vals = @(0, 1, 3, 3, 1, 3) },
    @{ Row = 4;  Name = "Usenix";    Vals = @(0, 0, 0, 1, 2, 2) },
    @{ Row = 5;  Name = "NDSS";      Vals = @(0, 0, 1, 1, 2, 3) },
    @{ Row = 6;  Name = "FC";        Vals = @(2, 3, 4, 7, 7, 5) },
    @{ Row = 7;  Name = "Crypto";    Vals = @(0, 0, 0, 2, 0, 4) },
    @{ Row = 8;  Name = "Eurocrypt"; Vals = @(0, 0, 0, 2, 1, 2) },
    @{ Row = 9;  Name = "Asiacrypt"; Vals = @(0, 0, 0, 0, 1, 2) },
    @{ Row = 10; Name = "TCC";       Vals = @(0, 0, 0, 0, 2, 3) }
)

foreach ($conf in $confRows) {
    $r = $conf.Row
    $byConfYear.Cells.Item($r, 1).Value = $conf.Name
    for ($i = 0; $i -lt $conf.Vals.Length; $i++) {
        $byConfYear.Cells.Item($r, 2 + $i).Value = $conf.Vals[$i]
    }
}

# Row total (H column): plain formula on row 2, shared formula H3:H11.
$byConfYear.Range("H2").Formula = "=SUM(B2:G2)"
$byConfYear.Range("H3:H11").Formula = "=SUM(B3:G3)"

# Totals row (row 11): "Totals" label (bold) + column sums B:G.
$byConfYear.Range("A11").Value = "Totals"
$byConfYear.Range("A11").Font.Bold = $true
$byConfYear.Range("B11").Formula = "=SUM(B2:B10)"
$byConfYear.Range("C11:G11").Formula = "=SUM(C2:C10)"

# ---------------------------------------------------------------------
# 3) Selection / active-sheet state to match the saved workbook.
# ---------------------------------------------------------------------
[void]$byTopic.Range("G13").Select()
[void]$byConfYear.Range("F6").Select()
